# Lecture 4 "Real world data" deck update
# - Insert a new "Lecture 4 outline" slide as slide 2 (everything else shifts down)
# - Update slide 1 subtitle text
# - Refresh the date placeholder text (28/11/2019 -> 02/12/2019) on every layout + the master
# - Switch the slide master background from a flat bgRef to a gradient fill

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 1: subtitle text change
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Real world data sources"

# ---------------------------------------------------------------------------
# 2. Insert a brand new slide at position 2 ("Lecture 4 outline")
# ---------------------------------------------------------------------------
$outline = $p.Slides.Add(2, 2)

$outline.Shapes.Item(1).TextFrame.TextRange.Text = "Lecture 4 outline"

$contentTr = $outline.Shapes.Item(2).TextFrame.TextRange
$contentTr.Text = "The importance of data (~10mins)`r" + `
  "Emergencies (~10mins)`r" + `
  "Emergency example (~10mins)`r" + `
  "Data gathering vs data generating: experiments (~10mins)`r" + `
  "Data gathering vs data generating: real-world sources (~10mins)`r"
$contentTr.Paragraphs(5).Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Refresh the "datetimeFigureOut" placeholder text everywhere
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
  $layout = $master.CustomLayouts.Item($i)
  for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
    $shp = $layout.Shapes.Item($j)
    if ($shp.Name -like "Date Placeholder*") {
      $shp.TextFrame.TextRange.Text = "02/12/2019"
    }
  }
}

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
  $shp = $master.Shapes.Item($j)
  if ($shp.Name -like "Date Placeholder*") {
    $shp.TextFrame.TextRange.Text = "02/12/2019"
  }
}

# ---------------------------------------------------------------------------
# 4. Slide master background: solid bgRef -> gradient fill
# ---------------------------------------------------------------------------
$master.Background.Fill.TwoColorGradient(1, 1)

Write-Output "edit complete"
